# Cell_Library.xlsx — "changing layout for the EIAN analysis workbook"
#
# On the "E_I and A_N ratios" sheet, the single "cells in pair" column
# (values like "1", "2" or "1,2") is replaced by two boolean columns,
# "analyze cell 1" and "analyze cell 2", inserted right after the
# "pair/cell number" column. Everything that used to live in columns
# D..H shifts one column to the right (E..I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E_I and A_N ratios")

# Insert a new blank column at D; old column C ("cells in pair") stays put,
# old columns D..H (area, type cell 1, type cell 2, layer, comments) shift to E..I.
$ws.Columns.Item(4).Insert()

# Recompute columns C & D (analyze cell 1 / analyze cell 2) from the old
# "cells in pair" values that are still sitting in column C, then stamp
# the new header labels.
for ($r = 2; $r -le 15; $r++) {
    $oldVal = [string]$ws.Cells.Item($r, 3).Value()
    $hasCell1 = $oldVal.Contains("1")
    $hasCell2 = $oldVal.Contains("2")

    if ($hasCell1) { $ws.Cells.Item($r, 3).Value = 1 } else { $ws.Cells.Item($r, 3).Value = 0 }
    if ($hasCell2) { $ws.Cells.Item($r, 4).Value = 1 } else { $ws.Cells.Item($r, 4).Value = 0 }
}

# Set D1 before C1 so the shared-string table gets "analyze cell 2" (76)
# ahead of "analyze cell 1" (77), matching the canonical save order.
$ws.Cells.Item(1, 4).Value = "analyze cell 2"
$ws.Cells.Item(1, 3).Value = "analyze cell 1"

# Leave the sheet's selection parked on the new header cell, as in the
# saved file.
$ws.Activate() | Out-Null
$ws.Range("C1").Select() | Out-Null
